# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" positioned between "2021-Q4" and "总计",
#    populated with the fund-holdings detail rows for the quarter.
# 2. Insert a new summary row at the top of "总计" ("2022-Q1", 9, 0.67),
#    pushing the existing rows down.

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------------
# 1) New "2022-Q1" sheet, inserted right after "2021-Q4" (i.e. right before 总计)
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $q4)
$newSheet.Name = "2022-Q1"

# NB: fetch "总计" only *after* the sheet insertion above - grabbing it
# beforehand would leave a reference that now resolves to the freshly
# inserted sheet once the collection shifts.
$zj = $wb.Worksheets.Item("总计")

# Match the page margins used by the other data sheets.
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Borrow the header/body cell formatting (fonts, borders, alignment) from the
# "2021-Q4" sheet so the new sheet matches the rest of the workbook: row 1 is
# the bold/bordered header style, and column A uses the same style down every
# data row.
$q4.Range("A1:H2").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)
$newSheet.Range("A2").Copy()
$newSheet.Range("A3:A10").PasteSpecial(-4122)

$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "'001208"
$newSheet.Cells.Item(2,3).Value = "诺安低碳经济股票A"
$newSheet.Cells.Item(2,4).Value = "'14.51"
$newSheet.Cells.Item(2,5).Value = "'82.03"
$newSheet.Cells.Item(2,6).Value = "'1.70"
$newSheet.Cells.Item(2,7).Value = "'0.2467"
$newSheet.Cells.Item(2,8).Value = 7

$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "'012071"
$newSheet.Cells.Item(3,3).Value = "中加喜利回报一年持有期混合A"
$newSheet.Cells.Item(3,4).Value = "'5.20"
$newSheet.Cells.Item(3,5).Value = "'46.57"
$newSheet.Cells.Item(3,6).Value = "'1.96"
$newSheet.Cells.Item(3,7).Value = "'0.1019"
$newSheet.Cells.Item(3,8).Value = 7

$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).Value = "'005775"
$newSheet.Cells.Item(4,3).Value = "中加转型动力灵活配置混合A"
$newSheet.Cells.Item(4,4).Value = "'3.41"
$newSheet.Cells.Item(4,5).Value = "'66.34"
$newSheet.Cells.Item(4,6).Value = "'2.33"
$newSheet.Cells.Item(4,7).Value = "'0.0795"
$newSheet.Cells.Item(4,8).Value = 9

$newSheet.Cells.Item(5,1).Value = 3
$newSheet.Cells.Item(5,2).Value = "'009242"
$newSheet.Cells.Item(5,3).Value = "中加核心智造混合A"
$newSheet.Cells.Item(5,4).Value = "'2.05"
$newSheet.Cells.Item(5,5).Value = "'65.71"
$newSheet.Cells.Item(5,6).Value = "'2.98"
$newSheet.Cells.Item(5,7).Value = "'0.0611"
$newSheet.Cells.Item(5,8).Value = 6

$newSheet.Cells.Item(6,1).Value = 4
$newSheet.Cells.Item(6,2).Value = "'010349"
$newSheet.Cells.Item(6,3).Value = "诺安低碳经济股票C"
$newSheet.Cells.Item(6,4).Value = "'3.52"
$newSheet.Cells.Item(6,5).Value = "'82.03"
$newSheet.Cells.Item(6,6).Value = "'1.70"
$newSheet.Cells.Item(6,7).Value = "'0.0598"
$newSheet.Cells.Item(6,8).Value = 7

$newSheet.Cells.Item(7,1).Value = 5
$newSheet.Cells.Item(7,2).Value = "'012072"
$newSheet.Cells.Item(7,3).Value = "中加喜利回报一年持有期混合C"
$newSheet.Cells.Item(7,4).Value = "'2.74"
$newSheet.Cells.Item(7,5).Value = "'46.57"
$newSheet.Cells.Item(7,6).Value = "'1.96"
$newSheet.Cells.Item(7,7).Value = "'0.0537"
$newSheet.Cells.Item(7,8).Value = 7

$newSheet.Cells.Item(8,1).Value = 6
$newSheet.Cells.Item(8,2).Value = "'005776"
$newSheet.Cells.Item(8,3).Value = "中加转型动力灵活配置混合C"
$newSheet.Cells.Item(8,4).Value = "'1.92"
$newSheet.Cells.Item(8,5).Value = "'66.34"
$newSheet.Cells.Item(8,6).Value = "'2.33"
$newSheet.Cells.Item(8,7).Value = "'0.0447"
$newSheet.Cells.Item(8,8).Value = 9

$newSheet.Cells.Item(9,1).Value = 7
$newSheet.Cells.Item(9,2).Value = "'320020"
$newSheet.Cells.Item(9,3).Value = "诺安策略精选股票"
$newSheet.Cells.Item(9,4).Value = "'1.33"
$newSheet.Cells.Item(9,5).Value = "'80.31"
$newSheet.Cells.Item(9,6).Value = "'1.73"
$newSheet.Cells.Item(9,7).Value = "'0.0230"
$newSheet.Cells.Item(9,8).Value = 8

$newSheet.Cells.Item(10,1).Value = 8
$newSheet.Cells.Item(10,2).Value = "'009243"
$newSheet.Cells.Item(10,3).Value = "中加核心智造混合C"
$newSheet.Cells.Item(10,4).Value = "'0.10"
$newSheet.Cells.Item(10,5).Value = "'65.71"
$newSheet.Cells.Item(10,6).Value = "'2.98"
$newSheet.Cells.Item(10,7).Value = "'0.0030"
$newSheet.Cells.Item(10,8).Value = 6

# Cells that carried a leading apostrophe (to force text / keep leading
# zeros) pick up an Excel "quote prefix" style; strip it so the cells end up
# with plain/default styling like the source data, while the values stay text.
$newSheet.Range("B2:G10").ClearFormats()

# ---------------------------------------------------------------------------
# 2) Insert a new top data row in "总计" for 2022-Q1, pushing the rest down.
# ---------------------------------------------------------------------------
$zj.Rows.Item(2).Insert()

$zj.Cells.Item(3,1).Copy()
$zj.Cells.Item(2,1).PasteSpecial(-4122)

$zj.Cells.Item(2,1).Value = 0
$zj.Cells.Item(2,2).Value = "2022-Q1"
$zj.Cells.Item(2,3).Value = 9
$zj.Cells.Item(2,4).Value = 0.67

$zj.Range("B2:D2").ClearFormats()

# The row-index column (A) holds a literal running index, not a formula, so
# the values that shifted down with the existing rows need to be renumbered
# to stay sequential (0, 1, 2, ...).
$zj.Cells.Item(3,1).Value = 1
$zj.Cells.Item(4,1).Value = 2

# Restore the originally-active sheet/selection (adding/editing sheets above
# shifts Excel's notion of the "active" tab to whatever we touched last).
$wb.Worksheets.Item("2021-Q3").Activate()
